$d = $word.ActiveDocument

$r = $d.Content
$r.Find.Execute("критического пути", $false, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)

$s = $r.Start
$e = $r.End

$target = $d.Range($s, $e)
$target.Text = "исправленного функционала"

$newLen = ("исправленного функционала").Length
$target2 = $d.Range($s, $s + $newLen)
# Toggle bold off/on so the engine keeps this as a distinct run from the
# preceding "тестирование " run instead of re-coalescing them.
$target2.Bold = 0
$target2.Bold = 1
